# This script reorders data rows 2-7 on the active worksheet.
# The underlying records (uniquely identified by column A, the "Id") are
# the same as before; only their row position within the table changes.
#
# Mapping of destination row -> source row (1-based worksheet rows):
#   2 <- 4
#   3 <- 7
#   4 <- 6
#   5 <- 2
#   6 <- 5
#   7 <- 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns that contain data for rows 2-7 (A through AY).
$firstCol = 1           # A
$lastCol  = 51           # AY

$sourceOfDest = @{
    2 = 4
    3 = 7
    4 = 6
    5 = 2
    6 = 5
    7 = 3
}

# 1) Snapshot the current values of every source row (cell by cell) before
#    overwriting anything, so row moves don't clobber each other.
$snapshot = @{}
foreach ($srcRow in 2..7) {
    $rowValues = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowValues[$col] = $ws.Cells.Item($srcRow, $col).Value2
    }
    $snapshot[$srcRow] = $rowValues
}

# 2) Write the snapshots back out to their destination rows, cell by cell.
#    Text values (e.g. date-like strings such as "2023-07-27") must be
#    written into a text-formatted cell, otherwise Excel's COM layer will
#    silently reinterpret them as real dates/numbers.
foreach ($destRow in 2..7) {
    $srcRow = $sourceOfDest[$destRow]
    $rowValues = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $value = $rowValues[$col]
        $destCell = $ws.Cells.Item($destRow, $col)
        if ($value -is [string]) {
            $destCell.NumberFormat = "@"
        }
        $destCell.Value = $value
    }
}
